$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Extend the "Project lifetime" header row (1) from year 10 out to year 15
#    by copying the format of L1 (style for header cells) into M1:Q1, then
#    filling in the sequential year numbers.
# ---------------------------------------------------------------------------
$ws.Range("L1").Copy()
$ws.Range("M1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("M1").Value = 11
$ws.Range("N1").Value = 12
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# ---------------------------------------------------------------------------
# 2. Initial Investment shrinks from -2,500,000 to -1,500,000
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = -1500000

# ---------------------------------------------------------------------------
# 3. Depreciation: extend C3:L3 pattern out through Q3, new flat value 30000
# ---------------------------------------------------------------------------
$ws.Range("L3").Copy()
$ws.Range("M3:Q3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3:Q3").Value = 30000

# ---------------------------------------------------------------------------
# 4. Incoming Payments: extend C4:L4 pattern out through Q4, new flat value 910000
# ---------------------------------------------------------------------------
$ws.Range("L4").Copy()
$ws.Range("M4:Q4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C4:Q4").Value = 910000

# ---------------------------------------------------------------------------
# 5. Outgoing Payments: extend C5:L5 pattern out through Q5, B5 & C5:Q5 change
# ---------------------------------------------------------------------------
$ws.Range("L5").Copy()
$ws.Range("M5:Q5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B5").Value = -140000
$ws.Range("C5:Q5").Value = -210000

# ---------------------------------------------------------------------------
# 6. Residual: value moves from L6 to Q6 (same 140000 figure)
# ---------------------------------------------------------------------------
$ws.Range("L6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("Q6").Value = 140000
$ws.Range("L6").Clear()

# ---------------------------------------------------------------------------
# 7. restricted Equity: B7 changes, and the offsetting figure moves from L7
#    to Q7 with a new (smaller) value
# ---------------------------------------------------------------------------
$ws.Range("L7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B7").Value = -100000
$ws.Range("Q7").Value = 100000
$ws.Range("L7").Clear()

# ---------------------------------------------------------------------------
# 8. Yearly Net: extend C8:L8 pattern out through Q8 (uses L8's style, which
#    already matches C8:K8's style), B8 changes, flat 730000 with a higher
#    final-year (Q8) value of 970000
# ---------------------------------------------------------------------------
$ws.Range("L8").Copy()
$ws.Range("M8:Q8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B8").Value = -1740000
$ws.Range("C8:P8").Value = 730000
$ws.Range("Q8").Value = 970000

# ---------------------------------------------------------------------------
# 9. Present Value: extend C9:L9 pattern out through Q9, and recompute every
#    year's discounted value for the new Yearly Net figures
# ---------------------------------------------------------------------------
$ws.Range("L9").Copy()
$ws.Range("M9:Q9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B9").Value = -1740000
$ws.Range("C9").Value = 673431.7343173431
$ws.Range("D9").Value = 621246.9873776227
$ws.Range("E9").Value = 573106.0769166261
$ws.Range("F9").Value = 528695.6429120167
$ws.Range("G9").Value = 487726.6078524138
$ws.Range("H9").Value = 449932.2950668024
$ws.Range("I9").Value = 415066.6928660539
$ws.Range("J9").Value = 382902.8531974667
$ws.Range("K9").Value = 353231.4143888069
$ws.Range("L9").Value = 325859.238366058
$ws.Range("M9").Value = 300608.1534742232
$ws.Range("N9").Value = 277313.7947179181
$ws.Range("O9").Value = 255824.5338726181
$ws.Range("P9").Value = 236000.4925024152
$ws.Range("Q9").Value = 289289.3870082176

# ---------------------------------------------------------------------------
# 10. Accumulated Present Value: extend C10:L10 pattern out through Q10 and
#     recompute the running totals. The red/green (negative/positive) style
#     boundary shifts from between H10/I10 to between D10/E10, since the
#     cumulative figure now turns positive a few years earlier.
# ---------------------------------------------------------------------------
$ws.Range("L10").Copy()
$ws.Range("M10:Q10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Cells E10:H10 flip from the "negative" (red) style to the "positive"
# (green) style used by I10:L10 - copy that format over first.
$ws.Range("I10").Copy()
$ws.Range("E10:H10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B10").Value = -1740000
$ws.Range("C10").Value = -1066568.265682657
$ws.Range("D10").Value = -445321.278305034
$ws.Range("E10").Value = 127784.798611592
$ws.Range("F10").Value = 656480.4415236088
$ws.Range("G10").Value = 1144207.049376023
$ws.Range("H10").Value = 1594139.344442825
$ws.Range("I10").Value = 2009206.037308879
$ws.Range("J10").Value = 2392108.890506345
$ws.Range("K10").Value = 2745340.304895152
$ws.Range("L10").Value = 3071199.54326121
$ws.Range("M10").Value = 3371807.696735434
$ws.Range("N10").Value = 3649121.491453352
$ws.Range("O10").Value = 3904946.02532597
$ws.Range("P10").Value = 4140946.517828385
$ws.Range("Q10").Value = 4430235.904836603

# ---------------------------------------------------------------------------
# 11. Net Present Value mirrors the final Accumulated Present Value figure
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = 4430235.904836603
